# Weekly forward look stats update 06.09
# The "as at" week has moved on: drop the week(s) that have now passed
# (rows for "02 Sep 2024" / "09 Sep 2024" / "16 Sep 2024") and refresh
# the "as at" date stamp near the top of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the introductory sentence with the new "as at" date.
$ws.Range("A2").Value = "This list contains a week-by-week view of  MoJ Official and National Statistics that have been pre-announced on the gov.uk release calendar as at 06 September 2024"

# Remove the three rows that have now passed (week commencing 02 Sep 2024,
# plus the two following blank weeks 09 Sep 2024 and 16 Sep 2024), shifting
# all subsequent weeks up so the table starts again at row 5.
$ws.Range("A5:F7").EntireRow.Delete()
